$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- Shift the September (R:S) column entries down by one row ---
# Old rows 35..90 move to new rows 36..91 (walk bottom-to-top so we never
# overwrite a source row before it has been read).
for ($r = 90; $r -ge 35; $r--) {
    $detail = $ws.Cells.Item($r, 18).Value()   # column R
    $date   = $ws.Cells.Item($r, 19).Value()   # column S
    $ws.Cells.Item($r + 1, 18).Value = $detail
    $ws.Cells.Item($r + 1, 19).Value = $date
}

# New entry at the top of the September list
$ws.Cells.Item(35, 18).Value = "saravanan"
$ws.Cells.Item(35, 19).Value = "2024-09-09 10:43:11"

# --- Shift the August (P:Q) column entries down by one row ---
# Old rows 91..94 move to new rows 92..95.
for ($r = 94; $r -ge 91; $r--) {
    $detail = $ws.Cells.Item($r, 16).Value()   # column P
    $date   = $ws.Cells.Item($r, 17).Value()   # column Q
    $ws.Cells.Item($r + 1, 16).Value = $detail
    $ws.Cells.Item($r + 1, 17).Value = $date
}

# Row 91 no longer carries an August entry
$ws.Cells.Item(91, 16).Value = ""
$ws.Cells.Item(91, 17).Value = ""

# --- The "Broadband" category label moves from row 95 to the new row 96 ---
$ws.Cells.Item(96, 1).Value = $ws.Cells.Item(95, 1).Value()
$ws.Cells.Item(95, 1).Value = ""
